# Scheduled runner: refresh Universalis market-price snapshots (current
# average prices) and the leve-profit figures derived from them, per
# crafting-job sheet (Leve Name/Item/Level/EXP/Gil/Amount/Item ID in A:G
# are left untouched).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2779.2666
$ws.Range("I17").Value = 3600
$ws.Range("J17").Value = 2480.818
$ws.Range("K17").Value = 10800
$ws.Range("L17").Value = 7442.454000000001
$ws.Range("M17").Value = -10632
$ws.Range("N17").Value = -7778.454000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 177.33333
$ws.Range("I33").Value = 179.11765
$ws.Range("J33").Value = 147
$ws.Range("K33").Value = 179.11765
$ws.Range("L33").Value = 147
$ws.Range("M33").Value = 49.88235

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2181.5
$ws.Range("I43").Value = 2737.5
$ws.Range("J43").Value = 1625.5
$ws.Range("K43").Value = 2737.5
$ws.Range("L43").Value = 1625.5
$ws.Range("M43").Value = -2668.5
$ws.Range("N43").Value = -1763.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 641.6667
$ws.Range("I127").Value = 518.2727
$ws.Range("J127").Value = 1999
$ws.Range("K127").Value = 1554.8181
$ws.Range("L127").Value = 5997
$ws.Range("M127").Value = 3405.1819

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 9400.725
$ws.Range("I132").Value = 7785.1055
$ws.Range("J132").Value = 12470.4
$ws.Range("K132").Value = 23355.3165
$ws.Range("L132").Value = 37411.2
$ws.Range("M132").Value = -20825.3165

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2811.718
$ws.Range("I138").Value = 4134.4116
$ws.Range("J138").Value = 2443.0984
$ws.Range("K138").Value = 12403.2348
$ws.Range("L138").Value = 7329.2952
$ws.Range("M138").Value = -7263.234800000002
$ws.Range("N138").Value = -17609.2952

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7056.577
$ws.Range("I32").Value = 7750.6665
$ws.Range("J32").Value = 6110.091
$ws.Range("K32").Value = 7750.6665
$ws.Range("L32").Value = 6110.091
$ws.Range("M32").Value = -7463.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2753.6296
$ws.Range("I61").Value = 1780.6842
$ws.Range("J61").Value = 5064.375
$ws.Range("K61").Value = 1780.6842
$ws.Range("L61").Value = 5064.375
$ws.Range("M61").Value = -1568.6842
$ws.Range("N61").Value = -5488.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2753.6296
$ws.Range("I136").Value = 1780.6842
$ws.Range("J136").Value = 5064.375
$ws.Range("K136").Value = 5342.0526
$ws.Range("L136").Value = 15193.125
$ws.Range("M136").Value = -2792.0526
$ws.Range("N136").Value = -20293.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 11340.857
$ws.Range("I99").Value = 4007.5
$ws.Range("J99").Value = 14274.2
$ws.Range("K99").Value = 4007.5
$ws.Range("L99").Value = 14274.2
$ws.Range("M99").Value = -2509.5
$ws.Range("N99").Value = -17270.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8388986
$ws.Range("I105").Value = 401755.2
$ws.Range("J105").Value = 41669116
$ws.Range("K105").Value = 401755.2
$ws.Range("L105").Value = 41669116
$ws.Range("M105").Value = -400008.2
$ws.Range("N105").Value = -41672610

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5352.3125
$ws.Range("I31").Value = 3893.7334
$ws.Range("J31").Value = 6639.294
$ws.Range("K31").Value = 3893.7334
$ws.Range("L31").Value = 6639.294
$ws.Range("M31").Value = -3598.7334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5352.3125
$ws.Range("I34").Value = 3893.7334
$ws.Range("J34").Value = 6639.294
$ws.Range("K34").Value = 3893.7334
$ws.Range("L34").Value = 6639.294
$ws.Range("M34").Value = -3691.7334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2722.5
$ws.Range("I58").Value = 1714.4546
$ws.Range("J58").Value = 3374.7646
$ws.Range("K58").Value = 1714.4546
$ws.Range("L58").Value = 3374.7646
$ws.Range("M58").Value = -1511.4546
$ws.Range("N58").Value = -3780.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3362.4736
$ws.Range("I134").Value = 3059.2
$ws.Range("J134").Value = 4499.75
$ws.Range("K134").Value = 9177.599999999999
$ws.Range("L134").Value = 13499.25
$ws.Range("M134").Value = -6642.599999999999
$ws.Range("N134").Value = -18569.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2722.5
$ws.Range("I136").Value = 1714.4546
$ws.Range("J136").Value = 3374.7646
$ws.Range("K136").Value = 5143.3638
$ws.Range("L136").Value = 10124.2938
$ws.Range("M136").Value = -2593.3638
$ws.Range("N136").Value = -15224.2938

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 7000
$ws.Range("I120").Value = 7000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 21000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -16162

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2399.2
$ws.Range("I132").Value = 1475.7142
$ws.Range("J132").Value = 2896.4614
$ws.Range("K132").Value = 13281.4278
$ws.Range("L132").Value = 26068.1526
$ws.Range("M132").Value = -10751.4278
$ws.Range("N132").Value = -31128.1526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4877.1904
$ws.Range("I122").Value = 3145
$ws.Range("J122").Value = 7692
$ws.Range("K122").Value = 9435
$ws.Range("L122").Value = 23076
$ws.Range("M122").Value = -6985
$ws.Range("N122").Value = -27976

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2719.5
$ws.Range("I126").Value = 2549.375
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 7648.125
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -5178.125
$ws.Range("N126").Value = -15140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 86999.5
$ws.Range("I137").Value = 90000
$ws.Range("J137").Value = 85999.336
$ws.Range("K137").Value = 90000
$ws.Range("L137").Value = 85999.336
$ws.Range("M137").Value = -84900
$ws.Range("N137").Value = -96199.336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 70677.39999999999
$ws.Range("I141").Value = 60390
$ws.Range("J141").Value = 73249.25
$ws.Range("K141").Value = 60390
$ws.Range("L141").Value = 73249.25
$ws.Range("M141").Value = -55210
$ws.Range("N141").Value = -83609.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 9155.691999999999
$ws.Range("I61").Value = 751.5833
$ws.Range("J61").Value = 110005
$ws.Range("K61").Value = 751.5833
$ws.Range("L61").Value = 110005
$ws.Range("M61").Value = -549.5833
$ws.Range("N61").Value = -110409

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 10000
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -10676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 10000
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -12340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 9155.691999999999
$ws.Range("I113").Value = 751.5833
$ws.Range("J113").Value = 110005
$ws.Range("K113").Value = 751.5833
$ws.Range("L113").Value = 110005
$ws.Range("M113").Value = 1418.4167
$ws.Range("N113").Value = -114345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 47676.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 47676.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 47676.5
$ws.Range("N54").Value = -48716.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6118
$ws.Range("I62").Value = 3884
$ws.Range("J62").Value = 6676.5
$ws.Range("K62").Value = 3884
$ws.Range("L62").Value = 6676.5
$ws.Range("M62").Value = -3260
$ws.Range("N62").Value = -7924.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6118
$ws.Range("I65").Value = 3884
$ws.Range("J65").Value = 6676.5
$ws.Range("K65").Value = 19420
$ws.Range("L65").Value = 33382.5
$ws.Range("M65").Value = -16300
$ws.Range("N65").Value = -39622.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3295.3809
$ws.Range("I81").Value = 1755.3334
$ws.Range("J81").Value = 5348.778
$ws.Range("K81").Value = 3510.6668
$ws.Range("L81").Value = 10697.556
$ws.Range("M81").Value = -2449.6668
$ws.Range("N81").Value = -12819.556

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3295.3809
$ws.Range("I84").Value = 1755.3334
$ws.Range("J84").Value = 5348.778
$ws.Range("K84").Value = 17553.334
$ws.Range("L84").Value = 53487.78
$ws.Range("M84").Value = -12249.334
$ws.Range("N84").Value = -64095.78

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 76927980
$ws.Range("I136").Value = 125001280
$ws.Range("J136").Value = 10711
$ws.Range("K136").Value = 375003840
$ws.Range("L136").Value = 32133
$ws.Range("M136").Value = -375001290
$ws.Range("N136").Value = -37233
